# The workbook's first sheet ("Sheet1", the "Data_Provider" data sheet)
# contains a "Result" column (D) whose data rows (D2:D5) used to hold
# "Pass" stamped with a colored fill. The user selected that range and
# cleared it completely (contents + formatting), leaving only the D1
# header ("Result") behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rng = $ws.Range("D2:D5")
$rng.Select()
$rng.Clear()
